$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 17139.777
$ws.Range("I43").Value = 2257.2856
$ws.Range("J43").Value = 33167.08
$ws.Range("K43").Value = 2257.2856
$ws.Range("L43").Value = 33167.08
$ws.Range("M43").Value = -2188.2856
$ws.Range("N43").Value = -33305.08
$ws.Range("H112").Value = 2721.7058
$ws.Range("J112").Value = 2916.4482
$ws.Range("L112").Value = 8749.3446
$ws.Range("N112").Value = -10965.3446
$ws.Range("H138").Value = 2974.0195
$ws.Range("J138").Value = 3061.7805
$ws.Range("L138").Value = 9185.341499999999
$ws.Range("N138").Value = -19465.3415

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1242.25
$ws.Range("I2").Value = 1087.6129
$ws.Range("K2").Value = 1087.6129
$ws.Range("M2").Value = -974.6129000000001
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15976
$ws.Range("H55").Value = 33142.715
$ws.Range("H61").Value = 4395
$ws.Range("I61").Value = 4136.769
$ws.Range("K61").Value = 4136.769
$ws.Range("M61").Value = -3924.769
$ws.Range("H80").Value = 53525
$ws.Range("J80").Value = 54666.668
$ws.Range("L80").Value = 54666.668
$ws.Range("N80").Value = -56662.668
$ws.Range("H83").Value = 53525
$ws.Range("J83").Value = 54666.668
$ws.Range("L83").Value = 164000.004
$ws.Range("N83").Value = -173984.004
$ws.Range("H97").Value = 1290.1724
$ws.Range("I97").Value = 1294.8636
$ws.Range("K97").Value = 1294.8636
$ws.Range("M97").Value = -798.8635999999999
$ws.Range("H116").Value = 1242.25
$ws.Range("I116").Value = 1087.6129
$ws.Range("K116").Value = 1087.6129
$ws.Range("M116").Value = 1206.3871
$ws.Range("H136").Value = 4395
$ws.Range("I136").Value = 4136.769
$ws.Range("K136").Value = 12410.307
$ws.Range("M136").Value = -9860.307000000001

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1242.25
$ws.Range("I3").Value = 1087.6129
$ws.Range("K3").Value = 1087.6129
$ws.Range("M3").Value = -973.6129000000001
$ws.Range("H12").Value = 7500
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -9832
$ws.Range("N12").Value = -5336
$ws.Range("H86").Value = 1617.6666
$ws.Range("I86").Value = 1672.0834
$ws.Range("K86").Value = 1672.0834
$ws.Range("M86").Value = -549.0834
$ws.Range("H89").Value = 1617.6666
$ws.Range("I89").Value = 1672.0834
$ws.Range("K89").Value = 8360.416999999999
$ws.Range("M89").Value = -2744.416999999999
$ws.Range("H99").Value = 1131.7778
$ws.Range("I99").Value = 1038.8
$ws.Range("K99").Value = 1038.8
$ws.Range("M99").Value = 459.2
$ws.Range("H107").Value = 1473.7
$ws.Range("I107").Value = 1473.7
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1473.7
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 446.3
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 68748
$ws.Range("J122").Value = 68748
$ws.Range("L122").Value = 68748
$ws.Range("N122").Value = -78548
$ws.Range("H134").Value = 4239
$ws.Range("I134").Value = 3186.8
$ws.Range("K134").Value = 9560.400000000001
$ws.Range("M134").Value = -7025.400000000001

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 586
$ws.Range("I13").Value = 724.5
$ws.Range("J13").Value = 447.5
$ws.Range("K13").Value = 724.5
$ws.Range("L13").Value = 447.5
$ws.Range("M13").Value = -585.5
$ws.Range("N13").Value = -725.5
$ws.Range("H86").Value = 6422.5713
$ws.Range("J86").Value = 6651.8
$ws.Range("L86").Value = 6651.8
$ws.Range("N86").Value = -8897.799999999999
$ws.Range("H89").Value = 6422.5713
$ws.Range("J89").Value = 6651.8
$ws.Range("L89").Value = 33259
$ws.Range("N89").Value = -44491
$ws.Range("H107").Value = 948.1429000000001
$ws.Range("J107").Value = 1081.5714
$ws.Range("L107").Value = 1081.5714
$ws.Range("N107").Value = -4921.5714
$ws.Range("H134").Value = 2641.111
$ws.Range("I134").Value = 2641.111
$ws.Range("K134").Value = 7923.333
$ws.Range("M134").Value = -5388.333

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 378.875
$ws.Range("I8").Value = 378.875
$ws.Range("K8").Value = 1136.625
$ws.Range("M8").Value = -997.625
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1163.6471
$ws.Range("I97").Value = 1125.7333
$ws.Range("K97").Value = 1125.7333
$ws.Range("M97").Value = -629.7333000000001
$ws.Range("H113").Value = 2352.1
$ws.Range("I113").Value = 2384.7144
$ws.Range("K113").Value = 2384.7144
$ws.Range("M113").Value = -214.7143999999998
$ws.Range("H132").Value = 2510.3333
$ws.Range("I132").Value = 2792.6875
$ws.Range("K132").Value = 8378.0625
$ws.Range("M132").Value = -5848.0625
$ws.Range("H134").Value = 110000
$ws.Range("J134").Value = 110000
$ws.Range("L134").Value = 330000
$ws.Range("N134").Value = -335070
$ws.Range("H136").Value = 14500.637
$ws.Range("J136").Value = 14500.637
$ws.Range("L136").Value = 43501.911
$ws.Range("N136").Value = -48601.911

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1236.15
$ws.Range("I61").Value = 1117.1578
$ws.Range("K61").Value = 1117.1578
$ws.Range("M61").Value = -915.1578
$ws.Range("H100").Value = 1976.6666
$ws.Range("I100").Value = 2020.25
$ws.Range("K100").Value = 2020.25
$ws.Range("M100").Value = -1479.25
$ws.Range("H113").Value = 1236.15
$ws.Range("I113").Value = 1117.1578
$ws.Range("K113").Value = 1117.1578
$ws.Range("M113").Value = 1052.8422
$ws.Range("H132").Value = 2961.3333
$ws.Range("I132").Value = 3108
$ws.Range("J132").Value = 2856.5715
$ws.Range("K132").Value = 9324
$ws.Range("L132").Value = 8569.7145
$ws.Range("M132").Value = -6794
$ws.Range("N132").Value = -13629.7145

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 895.25
$ws.Range("I107").Value = 593.7857
$ws.Range("J107").Value = 3005.5
$ws.Range("K107").Value = 1781.3571
$ws.Range("L107").Value = 9016.5
$ws.Range("M107").Value = 138.6428999999998
$ws.Range("N107").Value = -12856.5
$ws.Range("H126").Value = 3069.3572
$ws.Range("I126").Value = 2247.6667
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 6743.000100000001
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -4273.000100000001
$ws.Range("N126").Value = -28938.5
$ws.Range("H132").Value = 1378.8918
$ws.Range("I132").Value = 1369.1428
$ws.Range("K132").Value = 4107.428400000001
$ws.Range("M132").Value = -1577.428400000001
